$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data (row 2), pushing existing
# data rows down by 4. The sheet's used range stays fixed at A1:C21,
# so the last 4 existing data rows fall off the bottom (consistent
# with a rolling window of sensor readings).
$ws.Rows("2:5").Insert()

# Insert() copies formatting down from the row above (the bold header),
# so reset the newly inserted rows back to the plain/default style used
# by the rest of the data rows.
$ws.Range("A2:C5").Style = "Normal"
$ws.Range("A2:C5").ClearFormats()

# Populate the newly inserted rows with the new sensor readings.
$newData = @(
    @(0.01418807215633853, 0.06712245657330498, -0.07486735071454727),
    @(-0.02585268907603775, -0.07008951618557867, -0.06299911678901748),
    @(0.046578474342823, -0.1310305893421173, -0.0218384321779012),
    @(0.0740674127425465, -0.09423323614256736, -0.02838341776458984)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}

# Trim the sheet back down to its original extent (A1:C21) by removing
# the 4 rows that were pushed past the original last row.
$ws.Rows("22:25").Delete()
